# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Updates reconciliation names and recalculated Balance/Valeur Calculee/Jours de Stock
# for a handful of rows in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: LOUISE DAJEU
$ws.Range("B10").Value = "LOUISE DAJEU"
$ws.Range("F10").Value = 7518
$ws.Range("G10").Value = -5732
$ws.Range("H10").Value = 0.5673962264150944

# Row 11: ROSE DEUMENI
$ws.Range("B11").Value = "ROSE DEUMENI"
$ws.Range("F11").Value = 14870
$ws.Range("G11").Value = 5934.280000000001
$ws.Range("H11").Value = 1.664107648852023

# Row 13: Clarisse Ngenue Wankah
$ws.Range("B13").Value = "Clarisse Ngenue Wankah"
$ws.Range("F13").Value = 12048
$ws.Range("G13").Value = 6554.666666666667
$ws.Range("H13").Value = 2.193203883495146

# Row 15: Balance / Valeur Calculee / Jours de Stock update
$ws.Range("F15").Value = 35306
$ws.Range("G15").Value = 28086
$ws.Range("H15").Value = 4.890027700831025

# Row 17: ETIENNE JUSTIN JIOFACK
$ws.Range("B17").Value = "ETIENNE JUSTIN JIOFACK"
$ws.Range("F17").Value = 2717
$ws.Range("G17").Value = -4663.523076923077
$ws.Range("H17").Value = 0.3681310893119937

# Row 18: CLARISSE MAKOLO
$ws.Range("B18").Value = "CLARISSE MAKOLO"
$ws.Range("F18").Value = 7758
$ws.Range("G18").Value = -18661.66666666666
$ws.Range("H18").Value = 0.2936448857542993

# Row 19: LANDRY MANFOUO
$ws.Range("B19").Value = "LANDRY MANFOUO"
$ws.Range("F19").Value = 3855
$ws.Range("G19").Value = -35337.5
$ws.Range("H19").Value = 0.09836065573770492
